# Re-do TODE grade norms with week coding for grade
# Updates lookup B-column (scaled score) values for each grade/semester tab
# to reflect the new week-coded raw-score-to-scaled-score conversion.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("K-Fall")
$map = @{
    2 = 68
    3 = 71
    4 = 74
    5 = 77
    6 = 80
    7 = 83
    8 = 86
    9 = 89
    10 = 92
    11 = 95
    12 = 98
    13 = 101
    14 = 105
    15 = 108
    16 = 111
    17 = 114
    18 = 117
    19 = 120
    20 = 123
    21 = 126
    22 = 129
    23 = 130
}
foreach ($row in $map.Keys) {
    $ws.Cells.Item($row, 2).Value = $map[$row]
}

$ws = $wb.Worksheets.Item("K-Spring")
$map = @{
    2 = 61
    3 = 65
    4 = 68
    5 = 71
    6 = 74
    7 = 77
    8 = 80
    9 = 83
    10 = 86
    11 = 89
    14 = 98
    15 = 101
    16 = 104
    17 = 108
    18 = 111
    19 = 114
    20 = 117
    21 = 120
    22 = 123
    23 = 126
    24 = 129
}
foreach ($row in $map.Keys) {
    $ws.Cells.Item($row, 2).Value = $map[$row]
}

$ws = $wb.Worksheets.Item("1-Fall")
$map = @{
    2 = 55
    3 = 58
    4 = 61
    5 = 65
    6 = 68
    7 = 71
    8 = 74
    9 = 77
    10 = 80
    11 = 83
    13 = 89
    14 = 92
    15 = 95
    16 = 98
    17 = 101
    18 = 104
    19 = 107
    20 = 111
    21 = 114
    22 = 117
    23 = 120
    25 = 126
    26 = 129
}
foreach ($row in $map.Keys) {
    $ws.Cells.Item($row, 2).Value = $map[$row]
}

$ws = $wb.Worksheets.Item("1-Spring")
$map = @{
    2 = 68
    3 = 71
    4 = 74
    5 = 77
    6 = 80
    7 = 83
    8 = 86
    9 = 89
    10 = 92
    11 = 95
    12 = 98
    13 = 101
    14 = 105
    15 = 108
    16 = 111
    17 = 114
    18 = 117
    19 = 120
    20 = 123
    21 = 126
    22 = 129
    23 = 130
    24 = 130
    25 = 130
}
foreach ($row in $map.Keys) {
    $ws.Cells.Item($row, 2).Value = $map[$row]
}

$ws = $wb.Worksheets.Item("2-Fall")
$map = @{
    2 = 43
    3 = 46
    4 = 49
    5 = 52
    6 = 55
    7 = 58
    8 = 61
    9 = 64
    10 = 67
    12 = 74
    13 = 77
    14 = 80
    15 = 83
    16 = 86
    17 = 89
    18 = 92
    19 = 95
    20 = 98
    21 = 101
    22 = 104
    23 = 107
    24 = 110
    25 = 114
    26 = 117
}
foreach ($row in $map.Keys) {
    $ws.Cells.Item($row, 2).Value = $map[$row]
}

$ws = $wb.Worksheets.Item("2-Spring")
$map = @{
    2 = 40
    3 = 40
    4 = 43
    5 = 46
    6 = 49
    7 = 52
    8 = 55
    9 = 58
    10 = 61
    11 = 64
    12 = 67
    13 = 70
    14 = 74
    15 = 77
    16 = 80
    17 = 83
    18 = 86
    19 = 89
    20 = 92
    21 = 95
    23 = 101
    24 = 104
    25 = 107
    26 = 110
}
foreach ($row in $map.Keys) {
    $ws.Cells.Item($row, 2).Value = $map[$row]
}
